# newAssignment return from server
# Updates the per-worker assignment columns (C:G) with a fresh batch of
# values returned by the server: job titles in C/D/E, a numeric field in
# F, and a comma-separated ranking string in G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 1;  C = "Doctor";   D = "Teacher";  E = "";         F = 17; G = "4,5,1,3,6" }
    @{ Row = 2;  C = "Engineer"; D = "Doctor";   E = "Nurse";     F = 8;  G = "3,1,6,2" }
    @{ Row = 3;  C = "Engineer"; D = "Nurse";    E = "Doctor";    F = 11; G = "3,5,6,4" }
    @{ Row = 4;  C = "Teacher";  D = "Doctor";   E = "Nurse";     F = 20; G = "1,3,4,2,5" }
    @{ Row = 5;  C = "Nurse";    D = "Engineer"; E = "Doctor";    F = 19; G = "1,5,3,2,6" }
    @{ Row = 6;  C = "Doctor";   D = "Teacher";  E = "Nurse";     F = 14; G = "6,1,2,3" }
    @{ Row = 7;  C = "Doctor";   D = "Teacher";  E = "";          F = 20; G = "4,3,2,1,6,5" }
    @{ Row = 8;  C = "Doctor";   D = "Teacher";  E = "Engineer";  F = 13; G = "3,5,1,6,4" }
    @{ Row = 9;  C = "Engineer"; D = "Nurse";    E = "Doctor";    F = 18; G = "1,4,2,6,5,3" }
    @{ Row = 10; C = "Engineer"; D = "Nurse";    E = "Teacher";   F = 9;  G = "5,1,2,4,3,6" }
    @{ Row = 11; C = "Teacher";  D = "Nurse";    E = "Doctor";    F = 5;  G = "3,1,4,5,2,6" }
    @{ Row = 12; C = "Teacher";  D = "Engineer"; E = "Nurse";     F = 17; G = "4,2,6,5,1" }
    @{ Row = 13; C = "Teacher";  D = "Nurse";    E = "Doctor";    F = 8;  G = "3,5,2,4" }
    @{ Row = 14; C = "Doctor";   D = "Engineer"; E = "Nurse";     F = 13; G = "2,1,5,4,6,3" }
    @{ Row = 15; C = "Doctor";   D = "Nurse";    E = "";          F = 14; G = "3,1,2,6,5" }
    @{ Row = 16; C = "Doctor";   D = "Nurse";    E = "Engineer";  F = 14; G = "6,5,1,3,2" }
    @{ Row = 17; C = "Doctor";   D = "Teacher";  E = "Engineer";  F = 11; G = "3,4,1,6,5" }
    @{ Row = 18; C = "Engineer"; D = "Teacher";  E = "Doctor";    F = 10; G = "3,5,6,1,4" }
    @{ Row = 19; C = "Teacher";  D = "Doctor";   E = "";          F = 10; G = "2,3,1,6" }
    @{ Row = 20; C = "Nurse";    D = "Engineer"; E = "Doctor";    F = 16; G = "5,4,6,2,3" }
)

foreach ($rec in $rows) {
    $r = $rec.Row
    $ws.Cells.Item($r, 3).Value = $rec.C
    $ws.Cells.Item($r, 4).Value = $rec.D
    $ws.Cells.Item($r, 5).Value = $rec.E
    $ws.Cells.Item($r, 6).Value = $rec.F
    $ws.Cells.Item($r, 7).Value = $rec.G
}
